# Add two new vessel records (rows 11 and 12) to the active "2025-1" sheet,
# then move the current cell selection to G10, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11: REM TRITON (Remolcador / tug)
$ws.Range("E11").Value = "A.S/0032-125"
$ws.Range("D11").Value = "A.S/0032"
$ws.Range("B11").Value = "REM TRITON"
$ws.Range("A11").Value = "2025-1"
$ws.Range("C11").Value = "Remolcador"

# Row 12: EP MODESTO 8 (Embarcación Pesquera / fishing vessel)
$ws.Range("B12").Value = "EP MODESTO 8"
$ws.Range("A12").Value = "2025-1"
$ws.Range("C12").Value = "Embarcación Pesquera"
$ws.Range("E12").Value = "A.S/0033-125"
$ws.Range("D12").Value = "A.S/0033"

# Update the saved selection to G10, as in the source workbook.
[void]$ws.Range("G10").Select()
